$d = $word.ActiveDocument

# --- 1. Strip the emoji (Segoe UI Emoji / w16se:symEx "smiling face", U+1F60A) run ---
# The 5th paragraph in the document holds a single run whose content is an
# mc:AlternateContent (w16se:symEx + a "ðŸ˜Š" w:t fallback). That run is not
# exposed through the normal character-range navigation (its fallback text
# doesn't count towards Range.Start/End), so a plain Range.Delete() on the
# paragraph's (zero-length) range would just merge it away instead of
# clearing it. Work around this by first turning the phantom content into
# ordinary text (InsertAfter materializes/replaces it), then deleting that
# now-addressable text, which leaves a clean, empty paragraph behind.
$emojiPara = $d.Paragraphs(5).Range
$emojiPara.MoveEnd(1, -1)
$emojiPara.InsertAfter("X")

$emojiPara2 = $d.Paragraphs(5).Range
$emojiPara2.MoveEnd(1, -1)
$emojiPara2.Delete()

# --- 2. Remove the trailing paragraph that only contains the _GoBack bookmark ---
$d.Paragraphs.Last.Range.Delete()
